{"js": "// Update the date line and every two-digit multiplication answer in the\n// table, matching each old value to its replacement exactly once (in document\n// order), per the commit's regenerated \"answers\" sheet.\nconst replacements = [\n  [\"2024-10-05 Saturday\", \"2024-10-06 Sunday\"],\n  [\"31\u00d733=1023\", \"77\u00d797=7469\"],\n  [\"67\u00d739=2613\", \"13\u00d770=910\"],\n  [\"67\u00d791=6097\", \"36\u00d780=2880\"],\n  [\"43\u00d799=4257\", \"76\u00d770=5320\"],\n  [\"56\u00d740=2240\", \"23\u00d764=1472\"],\n  [\"17\u00d761=1037\", \"88\u00d763=5544\"],\n  [\"68\u00d782=5576\", \"42\u00d785=3570\"],\n  [\"79\u00d718=1422\", \"73\u00d775=5475\"],\n  [\"60\u00d730=1800\", \"70\u00d759=4130\"],\n  [\"82\u00d751=4182\", \"97\u00d741=3977\"],\n  [\"59\u00d744=2596\", \"32\u00d739=1248\"],\n  [\"48\u00d772=3456\", \"19\u00d724=456\"],\n  [\"32\u00d797=3104\", \"92\u00d763=5796\"],\n  [\"81\u00d779=6399\", \"15\u00d734=510\"],\n  [\"84\u00d738=3192\", \"67\u00d753=3551\"],\n  [\"97\u00d753=5141\", \"35\u00d755=1925\"],\n  [\"70\u00d726=1820\", \"22\u00d776=1672\"],\n  [\"36\u00d724=864\", \"81\u00d769=5589\"],\n  [\"87\u00d745=3915\", \"73\u00d785=6205\"],\n  [\"48\u00d784=4032\", \"14\u00d763=882\"],\n  [\"91\u00d741=3731\", \"37\u00d779=2923\"],\n  [\"45\u00d727=1215\", \"23\u00d725=575\"],\n  [\"49\u00d733=1617\", \"38\u00d723=874\"],\n  [\"14\u00d726=364\", \"90\u00d768=6120\"],\n  [\"18\u00d786=1548\", \"93\u00d736=3348\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Each old value is unique in the document, so replace the first (only) hit.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit multiplication answer in the\n# table, matching each old value to its replacement exactly once, per the\n# commit's regenerated \"answers\" sheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-05 Saturday\", \"2024-10-06 Sunday\"),\n    @(\"31\u00d733=1023\", \"77\u00d797=7469\"),\n    @(\"67\u00d739=2613\", \"13\u00d770=910\"),\n    @(\"67\u00d791=6097\", \"36\u00d780=2880\"),\n    @(\"43\u00d799=4257\", \"76\u00d770=5320\"),\n    @(\"56\u00d740=2240\", \"23\u00d764=1472\"),\n    @(\"17\u00d761=1037\", \"88\u00d763=5544\"),\n    @(\"68\u00d782=5576\", \"42\u00d785=3570\"),\n    @(\"79\u00d718=1422\", \"73\u00d775=5475\"),\n    @(\"60\u00d730=1800\", \"70\u00d759=4130\"),\n    @(\"82\u00d751=4182\", \"97\u00d741=3977\"),\n    @(\"59\u00d744=2596\", \"32\u00d739=1248\"),\n    @(\"48\u00d772=3456\", \"19\u00d724=456\"),\n    @(\"32\u00d797=3104\", \"92\u00d763=5796\"),\n    @(\"81\u00d779=6399\", \"15\u00d734=510\"),\n    @(\"84\u00d738=3192\", \"67\u00d753=3551\"),\n    @(\"97\u00d753=5141\", \"35\u00d755=1925\"),\n    @(\"70\u00d726=1820\", \"22\u00d776=1672\"),\n    @(\"36\u00d724=864\", \"81\u00d769=5589\"),\n    @(\"87\u00d745=3915\", \"73\u00d785=6205\"),\n    @(\"48\u00d784=4032\", \"14\u00d763=882\"),\n    @(\"91\u00d741=3731\", \"37\u00d779=2923\"),\n    @(\"45\u00d727=1215\", \"23\u00d725=575\"),\n    @(\"49\u00d733=1617\", \"38\u00d723=874\"),\n    @(\"14\u00d726=364\", \"90\u00d768=6120\"),\n    @(\"18\u00d786=1548\", \"93\u00d736=3348\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
